# aggiornamento fino a 13/03
# Adds 4 new daily rows (252-255) below the existing data, carrying the
# same formatting (date style) as the last existing data row (251),
# and extends the used data range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 251

$newDates = @(44326, 44327, 44328, 44329)
$newB     = @(0, 1, 0, 0)
$newC     = @(2, 3, 3, 2)
$newD     = @(67.43088334457181, 101.1463250168577, 101.1463250168577, 67.43088334457181)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $lastRow + 1 + $i

    # Copy the formatting (number format / style) of the last data row's
    # date cell so the new date cells keep the same look (style index).
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newB[$i]
    $ws.Cells.Item($r, 3).Value = $newC[$i]
    $ws.Cells.Item($r, 4).Value = $newD[$i]
}

$excel.CutCopyMode = $false
